# Deploying to main from @ pedalboard/pedalboard-soundcard@26f3ca50cc2e0b2273ee93242bc51fa4416f000b
# Update oscillator datasheet/supplier links (Abracon -> Kyocera/AVX parts) and
# refresh the "Created:" timestamp, plus increase row height on the two
# oscillator rows of the Costs sheet.

$wb = $excel.ActiveWorkbook

$wsBoM   = $wb.Worksheets.Item("BoM")
$wsCosts = $wb.Worksheets.Item("Costs")
$wsCostsDNF = $wb.Worksheets.Item("Costs (DNF)")

# --- Datasheet link: https://abracon.com/Oscillators/ASCO.pdf
#     -> https://media.digikey.com/pdf/Data%20Sheets/Kyocera%20International/Z_Series_X_Type.pdf
# Used on BoM rows 26 & 27 (column I) and Costs rows 27 & 28 (column E)
$newDatasheet = "https://media.digikey.com/pdf/Data%20Sheets/Kyocera%20International/Z_Series_X_Type.pdf"
$wsBoM.Range("I26").Value = $newDatasheet
$wsBoM.Range("I27").Value = $newDatasheet
$wsCosts.Range("E27").Value = $newDatasheet
$wsCosts.Range("E28").Value = $newDatasheet

# --- Supplier link for the 22.5792MHz oscillator (BoM row 26, column J)
$wsBoM.Range("J26").Value = "https://www.digikey.ch/en/products/detail/kyocera-avx/KC2016Z22-5792C1KX00/11610237"

# --- Supplier link for the 24.576MHz oscillator (BoM row 27, column J)
$wsBoM.Range("J27").Value = "https://www.digikey.ch/en/products/detail/kyocera-avx/KC2016Z24-5760C1KX00/11610181"

# --- "Created:" timestamp, shared between the Costs and Costs (DNF) sheets
$newCreated = "2024-12-15 20:10:23"
$wsCosts.Range("B31").Value = $newCreated
$wsCostsDNF.Range("B14").Value = $newCreated

# --- Row height bump on the oscillator rows of the Costs sheet
$wsCosts.Rows(27).RowHeight = 30
$wsCosts.Rows(28).RowHeight = 30
